$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The B1 cell ("AAAAA") carried a redundant explicit number-format style
# (General, but still flagged as "applied"). Drop it so the cell goes back
# to the default/unstyled state.
$ws.Range("B1").ClearFormats()

# New test data: row 8, used to exercise conditional formatting on a
# formatted number cell (formatRawCellContents() used to ignore CF).
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 12
$ws.Range("A8:B8").NumberFormat = "0.0"

# New conditional formatting rule: highlight A8:B8 when the value is
# greater than 10, using the same red-on-red style as the existing rules.
$rng = $ws.Range("A8:B8")
$fc = $rng.FormatConditions.Add(1, 5, "10")
$fc.Font.Color = 393372
$fc.Interior.Color = 13551615
$fc.SetFirstPriority()

# Move the active selection to B8, matching the saved workbook state.
$ws.Range("B8").Select()
